# Apply cryptos list update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.255.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.859.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7024"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08207"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3042"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08181"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7165"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.189"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.268.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007887"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.777"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.105.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.459"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.991"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1452"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.991"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.12%  "
$ws.Range("E30").Value = "  +4.44%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.486"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.407"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.053"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05220"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7077"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.39%  "
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01848"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.720"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.143.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9246"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.973"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.777"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.003.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.951"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.09%  "
